$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing customer's name (row 2)
$ws.Range("B2").Value = "Het B. Patel"

# Add new row 3
$ws.Range("A3").Value = "63e22be88db3f87bb229bb74"
$ws.Range("B3").Value = "Het B. Patel"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "230"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "KOT"

# Add new row 4
$ws.Range("A4").Value = "63e22d478db3f87bb229bc52"
$ws.Range("B4").Value = "Ayushi"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "10"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "10"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "KOT"
